# Add a new "range_detection" column (L) to the tracks_description sheet.
# For every data row this is "manual", except the Sonic curve row (row 7)
# which is set to "auto" - letting the curve scale be auto-detected from
# the 5th/95th percentile of the values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tracks_description")

# Header cell, formatted like the other header cells (copy format from the
# neighbouring "scale" header so borders/centering match).
$ws.Range("J1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null
$ws.Range("L1").Value = "range_detection"

# Data cells: default to "manual" for every row, matching the formatting
# used by the other text columns (copy format from column A).
$ws.Range("A2:A17").Copy() | Out-Null
$ws.Range("L2:L17").PasteSpecial(-4122) | Out-Null

# Row 7 is the Sonic curve - flip it to "auto" (written first so the shared
# string table gets "auto" before "manual", matching the source order).
$ws.Cells.Item(7, 12).Value = "auto"

foreach ($r in 2..17) {
    if ($r -ne 7) {
        $ws.Cells.Item($r, 12).Value = "manual"
    }
}

$ws.Application.CutCopyMode = $false

# Size the new column to fit its contents.
$ws.Columns.Item(12).ColumnWidth = 13.6

# Restore selection to match the editor's final cursor position.
$ws.Range("I20").Select() | Out-Null
